# Memory Benchmark: rename existing sheet, add a second sheet with heap-report
# test data, and switch focus to the new sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Rename the original sheet.
$ws1.Name = "VS Heap Tool"

# 2. Add the new "Heap Report from Test" sheet right after it.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Heap Report from Test"

# 3. Populate cells in the exact order the original author did, so that new
#    shared-string entries land at the same indices as in the source workbook.
$ws2.Range("D1").Value2 = "Branch"
$ws2.Range("F1").Value2 = "Heap Memory Test 1"
$ws2.Range("G1").Value2 = "Heap Memory Test 2"
$ws2.Range("H1").Value2 = "Heap Memory Test 3"
$ws2.Range("D2").Value2 = "Research_Memory_Management_Options"
$ws2.Range("B1").Value2 = "Platform"
$ws2.Range("B2").Value2 = "PC"
$ws2.Range("C1").Value2 = "Build"
$ws2.Range("C2").Value2 = "Release"
$ws2.Range("D3").Value2 = "Factory_Class"

# 4. Remaining header/text cells that reuse already-existing shared strings.
$ws2.Range("A1").Value2 = "Date"
$ws2.Range("E1").Value2 = "Duration"
$ws2.Range("I1").Value2 = "Description"
$ws2.Range("J1").Value2 = "Observations"
$ws2.Range("B3").Value2 = "PC"
$ws2.Range("C3").Value2 = "Release"

# 5. Data rows.
$ws2.Range("A2").Value2 = 43409.84375
$ws2.Range("A3").Value2 = 43409.854166666664

$ws2.Range("F2").Value2 = 255320430
$ws2.Range("G2").Value2 = 255320430
$ws2.Range("H2").Value2 = 255320430

$ws2.Range("F3").Value2 = 83737598
$ws2.Range("G3").Value2 = 190402070
$ws2.Range("H3").Value2 = 403696526

# 6. Copy number formats from the first sheet so style indices line up the
#    same way they do in the target workbook (date format + thousands/comma
#    format + wrap-text format).
$ws1.Range("A2").Copy()
$ws2.Range("A2:A3").PasteSpecial(-4122)

$ws1.Range("C1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("F1:H1").PasteSpecial(-4122)
$ws2.Range("F2:H3").PasteSpecial(-4122)

$ws1.Range("D1").Copy()
$ws2.Range("I1:J1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 7. Column widths on the new sheet (A, D and F:H are best-fit to their
#    contents; I:J are manually widened, mirroring sheet1's description /
#    observations columns).
$ws2.Columns.Item(1).EntireColumn.AutoFit()
$ws2.Columns.Item(4).EntireColumn.AutoFit()
$ws2.Range("F1:H1").EntireColumn.AutoFit()

$ws2.Columns.Item(1).ColumnWidth = 13.917
$ws2.Columns.Item(4).ColumnWidth = 35.2505
$ws2.Range("F1:H1").EntireColumn.ColumnWidth = 17.25
$ws2.Range("I1:J1").EntireColumn.ColumnWidth = 49.75

# 8. Sheet2 view: zoomed in, selection on H14.
$ws2.Range("H14").Select()
$excel.ActiveWindow.Zoom = 150

# 9. Sheet1 view: selection becomes the header row A1:E1.
$ws1.Range("A1:E1").Select()

# 10. Make the new sheet the active tab, matching the authored workbook.
$ws2.Activate()
